# Update the trial-2 row of the training schedule (x_corrSteps, x_nrSteps,
# and the derived praclen-ish total in column H), then leave the selection
# on D2 to match the author's last active cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 46

$ws.Range("D2").Select()
